# myopic/input/longtable.xlsx - "changed run loop, and myopic input creation"
#
# The upstream Python run loop that generates this workbook was changed, so a
# re-run produced a different CO2-limit input (row 2) and different
# commodity-price projections (rows 5-12) on Sheet1. Several of the
# commodity-price cells that used to carry a "grow 2%/year" formula now only
# carry the freshly-computed literal result (the generator script stopped
# emitting a formula for those years), while a couple of rows keep growing
# from a changed base value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: CO2 limit base value changes; F2/G2/H2 keep their "=prev*0.5"
#     formulas and simply recalculate off the new E2. ---
$ws.Range("E2").Value = 40

# --- Row 5 (Nuclear price): growth formulas replaced by their already
#     computed literal results (values unchanged by the regen). ---
$ws.Range("F5").Value = 2.6520000000000001
$ws.Range("G5").Value = 2.7050400000000003
$ws.Range("H5").Value = 2.7591408000000004
$ws.Range("I5").Value = 2.8143236160000002
$ws.Range("J5").Value = 2.8706100883200003
$ws.Range("K5").Value = 2.9280222900864006

# --- Row 6 (Biomass price): same - literal values only, unchanged numbers. ---
$ws.Range("F6").Value = 7.1400000000000006
$ws.Range("G6").Value = 7.2828000000000008
$ws.Range("H6").Value = 7.4284560000000006
$ws.Range("I6").Value = 7.5770251200000009
$ws.Range("J6").Value = 7.7285656224000014
$ws.Range("K6").Value = 7.883136934848002

# --- Row 7 (Coal price): literal values, and the regenerated numbers changed. ---
$ws.Range("F7").Value = 27.741091302211302
$ws.Range("G7").Value = 47.332182604422606
$ws.Range("H7").Value = 56.805584594594592
$ws.Range("I7").Value = 66.278986584766585
$ws.Range("J7").Value = 72.757548574938568
$ws.Range("K7").Value = 79.236110565110565

# --- Row 8 (Gas price): literal values; F8 unchanged, G8:K8 recomputed. ---
$ws.Range("F8").Value = 20.502000000000002
$ws.Range("G8").Value = 27.992005949288941
$ws.Range("H8").Value = 31.736446522217683
$ws.Range("I8").Value = 35.480887095146429
$ws.Range("J8").Value = 37.93953455115517
$ws.Range("K8").Value = 40.398182007163911

# --- Row 9 (Oil price): base value E9 changes, rest follow as literals. ---
$ws.Range("E9").Value = 50.46240705882353
$ws.Range("F9").Value = 50.46240705882353
$ws.Range("G9").Value = 50.46240705882353
$ws.Range("H9").Value = 56.764982647058829
$ws.Range("I9").Value = 63.067558235294122
$ws.Range("J9").Value = 67.009793823529421
$ws.Range("K9").Value = 70.952029411764713

# --- Row 10 (Other price): literal values, unchanged numbers. ---
$ws.Range("F10").Value = 56.1
$ws.Range("G10").Value = 57.222000000000001
$ws.Range("H10").Value = 58.366440000000004
$ws.Range("I10").Value = 59.533768800000004
$ws.Range("J10").Value = 60.724444176000006
$ws.Range("K10").Value = 61.938933059520004

# --- Row 11 (Elec sell price): base value E11 doubles, rest follow as literals. ---
$ws.Range("E11").Value = 20
$ws.Range("F11").Value = 20.399999999999999
$ws.Range("G11").Value = 20.808
$ws.Range("H11").Value = 21.224160000000001
$ws.Range("I11").Value = 21.648643200000002
$ws.Range("J11").Value = 22.081616064000002
$ws.Range("K11").Value = 22.523248385280002

# --- Row 12 (Elec buy price): formulas dropped, all years pinned to 100. ---
$ws.Range("F12").Value = 100
$ws.Range("G12").Value = 100
$ws.Range("H12").Value = 100
$ws.Range("I12").Value = 100
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 100

# --- Restore the editor's cursor/selection to F2 (was F6). ---
[void]$ws.Range("F2").Select()
